$d = $word.ActiveDocument

$pairs = @(
    @("32×27=864", "72×18=1296"),
    @("67×20=1340", "81×33=2673"),
    @("55×40=2200", "25×34=850"),
    @("44×59=2596", "50×43=2150"),
    @("46×57=2622", "24×21=504"),
    @("25×79=1975", "60×51=3060"),
    @("56×17=952", "46×38=1748"),
    @("91×63=5733", "54×87=4698"),
    @("44×75=3300", "22×18=396"),
    @("95×60=5700", "52×32=1664"),
    @("68×28=1904", "35×25=875"),
    @("85×14=1190", "58×29=1682"),
    @("64×49=3136", "72×48=3456"),
    @("82×24=1968", "17×55=935"),
    @("27×66=1782", "73×78=5694"),
    @("98×24=2352", "74×25=1850"),
    @("47×76=3572", "67×80=5360"),
    @("63×21=1323", "69×54=3726"),
    @("47×94=4418", "89×70=6230"),
    @("35×76=2660", "39×91=3549"),
    @("13×13=169", "75×20=1500"),
    @("30×19=570", "18×35=630"),
    @("73×29=2117", "31×54=1674"),
    @("25×24=600", "23×23=529"),
    @("39×51=1989", "29×27=783")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
